$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other")

# Gene symbols (column B) were entered first, top-to-bottom ...
$ws.Cells.Item(109, 2).Value = "CHD1L"
$ws.Cells.Item(110, 2).Value = "MYO3A"
$ws.Cells.Item(111, 2).Value = "MYO5B"

# ... then the matching Ensembl gene ids (column A) were filled in out of order
$ws.Cells.Item(110, 1).Value = "ENSG00000095777"
$ws.Cells.Item(111, 1).Value = "ENSG00000167306"
$ws.Cells.Item(109, 1).Value = "ENSG00000131778"

# ... and finally the shared PMID reference (column D), styled like the rows above it
$ws.Cells.Item(109, 4).Value = "PMID: 39037077"
$ws.Cells.Item(109, 4).Font.Color = 2171169

$ws.Cells.Item(110, 4).Value = "PMID: 39037077"
$ws.Cells.Item(110, 4).Font.Color = 2171169

$ws.Cells.Item(111, 4).Value = "PMID: 39037077"
$ws.Cells.Item(111, 4).Font.Color = 2171169

$ws.Range("D112").Select() | Out-Null
